$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 3 values
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = "1004944-53.2023.4.06.3804"

# Add new row 4
$ws.Range("A4").Value = "agravo3.pdf"
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = "0001547-36.2008.4.01.3814"

# Add new row 5
$ws.Range("A5").Value = "agravo4.pdf"
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = "5064477-98.2024.8.13.0702"
